# "Working with RF RX." -- add expected/actual comparison columns on both the
# raw-scope sheet and the decoded-bits sheet, and flag mismatches with
# conditional formatting.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "scope_0_2" -- raw waveform samples
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1"    -- decoded bit table

# ---------------------------------------------------------------------------
# 1) scope_0_2: add a "delta" column D with a few spot-check timing formulas.
# ---------------------------------------------------------------------------
$ws1.Range("D9").Formula = "=A9-A3"
$ws1.Range("D9").NumberFormat = "0.00000"

$ws1.Range("D10").Formula = "=A217-A10"
$ws1.Range("D10").NumberFormat = "0.00000"

$ws1.Range("D218").Formula = "=A237-A218"
$ws1.Range("D218").NumberFormat = "0.00000"

# Conditional formatting: flag any decoded value greater than 1 anywhere in
# column C, plus a second (higher-priority / later-added) rule scoped to C2.
# Created in this order so the dxf + priority bookkeeping matches Excel's
# own incremental numbering.
$throwAway1 = $ws1.Range("Z1").FormatConditions.Add(1, 5, "9999")
$throwAway1.Interior.Color = 255
$throwAway1.Delete()

$condAll = $ws1.Range("C1:C1048576").FormatConditions.Add(1, 5, "1")
$condAll.Font.Color = 393372
$condAll.Interior.Color = 13551615

$throwAway2 = $ws1.Range("Z2").FormatConditions.Add(1, 5, "9999")
$throwAway2.Font.Color = 393372
$throwAway2.Interior.Color = 13551615
$throwAway2.Delete()

$condC2 = $ws1.Range("C2").FormatConditions.Add(1, 5, "1")
$condC2.Font.Color = 393372
$condC2.Interior.Color = 13551615

# View: zoom to 100% and park the selection on L35 (no more scrolled-down
# top-left cell).
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100
$ws1.Range("L35").Select()

# ---------------------------------------------------------------------------
# 2) Sheet1 (decoded bits): rebuild as label | raw-count | decoded-bit |
#    expected-bit | mismatch, for the RF RX capture (20 address bits + 4
#    data bits), then two footnote rows.
# ---------------------------------------------------------------------------
$ws2.Activate()

$labels  = @("A0","A1","A2","A3","A4","A5","A6","A7","A8","A9","A10","A11","A12","A13","A14","A15","A16","A17","A18","A19","D0","D1","D2","D3")
$counts  = @(2059,684,682,2061,2061,684,2061,2062,2063,683,2062,683,2062,683,683,2062,2064,2060,684,683,2062,682,684,1540)
$expect  = @(1,0,0,1,1,0,1,1,1,0,1,0,1,0,0,1,1,1,0,0,1,0,0,1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 1).Value = $labels[$i]
    $ws2.Cells.Item($r, 2).Value = $counts[$i]
    $ws2.Cells.Item($r, 4).Value = $expect[$i]
}

$ws2.Range("C1").Formula = "=IF(B1>1000,1,0)"
$ws2.Range("C2:C24").Formula = "=IF(B2>1000,1,0)"

$ws2.Range("E1").Formula = "=D1-C1"
$ws2.Range("E2:E24").Formula = "=D2-C2"

$ws2.Rows.Item(25).ClearContents()

$ws2.Cells.Item(26, 2).Value = "1 count = 500ns"
$ws2.Cells.Item(27, 2).Value = "sync = 10.6ms"

# Conditional formatting: highlight any row where decoded bit != expected bit.
$condMismatch = $ws2.Range("E1:E1048576").FormatConditions.Add(1, 4, "0")
$condMismatch.Interior.Color = 255

$ws2.Range("B1:B24").Select()

$wb.Application.Calculate()
